# Insert a new row before row 2 (shifts existing data down by one row)
# and populate it with the new "Bairnsdale" exposure site entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Insert()

# The freshly inserted row inherits the header row's bold/centered formatting
# by default; clear it so the new row matches the plain style used by the
# other data rows.
$ws.Range("A2:D2").ClearFormats()

$ws.Range("A2").Value = "Bairnsdale"
$ws.Range("B2").Value = "V/Line train - Bairnsdale to Melbourne"
$ws.Range("C2").Value = "30/12/20 12:45pm-4:30pm"
$ws.Range("D2").Value = "Case caught train from Bairnsdale to Caulfield"
